$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 ("H 72"), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()
